# "Updated symbol list on Thu Dec 15 20:42:35 UTC 2022 with GitHub Actions"
#
# The source sheet stores every data cell (columns B-E) as literal text
# (OOXML inlineStr), even the numeric-looking "Price" column D. Assigning a
# plain numeric-looking string via Range.Value (e.g. "261.02") would make
# Excel auto-convert it to a real number - which both changes the stored
# cell type and introduces binary floating-point noise (261.02 ->
# 261.01999999999998). Prefixing every value with a leading apostrophe
# forces Excel's normal "treat as text" behaviour, so the text is stored
# verbatim (exact digits/trailing zeros preserved) just like the original
# cells.
function Set-Text($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price refresh only (no reordering) -----------------------------------
Set-Text "D2"  "261.02"      # BNB
Set-Text "D4"  "6.198"       # HuobiToken
Set-Text "D5"  "0.06111"     # Cronos
Set-Text "D6"  "6.740"       # KuCoinToken
Set-Text "D7"  "3.501"       # GateToken
Set-Text "D8"  "1.356"       # FTXToken
Set-Text "D9"  "0.7985"      # MXToken
Set-Text "D10" "0.1579"      # WazirX
Set-Text "D11" "0.08043"     # MandalaExchangeToken
Set-Text "D14" "0.09295"     # BitMartToken
Set-Text "D15" "3.902"       # MCDex
Set-Text "D16" "0.001706"    # BitForexToken
Set-Text "D17" "0.04829"     # CoinExToken

# Rows 18-24 rerank: "One" jumps from rank 24 up to rank 18, and
# TigerCash/BitKan/HotbitToken/NitroEx/LEO/BTSEToken each shift down one
# row, picking up refreshed prices/labels along the way.
Set-Text "B18" "One"
Set-Text "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-Text "D18" "0.0006161"
Set-Text "E18" "17OneONEWorstin24h"

Set-Text "B19" "TigerCash"
Set-Text "C19" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-Text "D19" "0.006199"
Set-Text "E19" "18TigerCashTCH"

Set-Text "B20" "BitKan"
Set-Text "C20" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-Text "D20" "0.001102"
Set-Text "E20" "19BitKanKAN"

Set-Text "B21" "HotbitToken"
Set-Text "C21" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-Text "D21" "0.003402"
Set-Text "E21" "20HotbitTokenHTB"

Set-Text "B22" "NitroEx"
Set-Text "C22" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-Text "D22" "0.0001501"
Set-Text "E22" "21NitroExNTX"

Set-Text "B23" "LEO"
Set-Text "C23" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-Text "D23" "3.690"
Set-Text "E23" "22LEOLEO"

Set-Text "B24" "BTSEToken"
Set-Text "C24" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-Text "D24" "2.262"
Set-Text "E24" "23BTSETokenBTSE"

# More price-only refreshes ---------------------------------------------
Set-Text "D25" "0.3359"      # BitpandaEcosystemToken
Set-Text "D27" "0.0006166"   # UpBots
Set-Text "D40" "0.04594"     # IDEX

# Rows 41-43: BKEXToken and KickToken swap ranks (CEJI stays put),
# along with refreshed prices.
Set-Text "B41" "BKEXToken"
Set-Text "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-Text "D41" "0.1124"
Set-Text "E41" "40BKEXTokenBKK"

Set-Text "D42" "0.003132"    # CEJI

Set-Text "B43" "KickToken"
Set-Text "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-Text "D43" "0.003414"
Set-Text "E43" "42KickTokenKICK"

# Trailing price-only refreshes ------------------------------------------
Set-Text "D46" "0.00006039"  # CoinLion
Set-Text "D48" "0.7501"      # CoinbaseStockToken
Set-Text "D49" "0.1185"      # BOLO (price jump, no longer "Worst in 24h")
Set-Text "E49" "48BOLOBOLO"
Set-Text "D50" "0.00001500"  # CryptobidCoin
